{"js": "// Word JS API (Office.js) script.\n// Runs as the body of: async (context) => { ... }\n//\n// Applies the three run-level edits described by the diff:\n//  1) Split \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\" so the run\n//     boundary moves from after \"NG\u01af\u1edcI \" to after \"NG\u01af\u1edcI \u0110\u1ea0I\", and the\n//     first run drops its custom rFonts/spacing (keeping bold/size/lang).\n//  2) Split \"{dissolution_approve_representative | upper}\" into three\n//     runs and rename the field to \"dissolution_approve_org_person\".\n\nconst OOXML_NS =\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPkg(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document ' + OOXML_NS + '><w:body>' +\n    bodyXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// --- Edit 1: \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\" heading ---\nconst headingResults = context.document.body.search(\n  \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\",\n  { matchCase: true, matchWholeWord: false }\n);\nheadingResults.load(\"items\");\nawait context.sync();\n\nif (headingResults.items.length === 0) {\n  throw new Error('Could not find heading text \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\"');\n}\n\nconst headingRange = headingResults.items[0];\nconst headingOoxml = wrapPkg(\n  \"<w:p>\" +\n    \"<w:r>\" +\n    \"<w:rPr><w:b/><w:sz w:val=\\\"26\\\"/><w:szCs w:val=\\\"26\\\"/><w:lang w:val=\\\"pt-BR\\\"/></w:rPr>\" +\n    \"<w:t>NG\u01af\u1edcI \u0110\u1ea0I</w:t>\" +\n    \"</w:r>\" +\n    \"<w:r>\" +\n    \"<w:rPr><w:b/><w:sz w:val=\\\"26\\\"/><w:szCs w:val=\\\"26\\\"/><w:lang w:val=\\\"pt-BR\\\"/></w:rPr>\" +\n    \"<w:t xml:space=\\\"preserve\\\"> DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P</w:t>\" +\n    \"</w:r>\" +\n    \"</w:p>\"\n);\nheadingRange.insertOoxml(headingOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: \"{dissolution_approve_representative | upper}\" placeholder ---\nconst placeholderResults = context.document.body.search(\n  \"{dissolution_approve_representative | upper}\",\n  { matchCase: true, matchWholeWord: false }\n);\nplaceholderResults.load(\"items\");\nawait context.sync();\n\nif (placeholderResults.items.length === 0) {\n  throw new Error('Could not find placeholder text \"{dissolution_approve_representative | upper}\"');\n}\n\nconst placeholderRange = placeholderResults.items[0];\nconst placeholderOoxml = wrapPkg(\n  \"<w:p>\" +\n    \"<w:r>\" +\n    \"<w:rPr><w:sz w:val=\\\"26\\\"/><w:szCs w:val=\\\"26\\\"/><w:lang w:val=\\\"pt-BR\\\"/></w:rPr>\" +\n    \"<w:t>{</w:t>\" +\n    \"</w:r>\" +\n    \"<w:r>\" +\n    \"<w:rPr><w:sz w:val=\\\"26\\\"/><w:szCs w:val=\\\"26\\\"/><w:lang w:val=\\\"pt-BR\\\"/></w:rPr>\" +\n    \"<w:t xml:space=\\\"preserve\\\">dissolution_approve_org_person </w:t>\" +\n    \"</w:r>\" +\n    \"<w:r>\" +\n    \"<w:rPr><w:sz w:val=\\\"26\\\"/><w:szCs w:val=\\\"26\\\"/><w:lang w:val=\\\"pt-BR\\\"/></w:rPr>\" +\n    \"<w:t>| upper}</w:t>\" +\n    \"</w:r>\" +\n    \"</w:p>\"\n);\nplaceholderRange.insertOoxml(placeholderOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $app / $doc are pre-seeded; the live document is $word.ActiveDocument.\n#\n# Applies the three run-level edits described by the diff:\n#  1) Split \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\" so the run\n#     boundary moves from after \"NG\u01af\u1edcI \" to after \"NG\u01af\u1edcI \u0110\u1ea0I\", and the\n#     first run drops its custom rFonts/spacing (keeping bold/size/lang).\n#  2) Split \"{dissolution_approve_representative | upper}\" into three\n#     runs and rename the field to \"dissolution_approve_org_person\".\n\n$d = $word.ActiveDocument\n\n$pkgOpen = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# --- Edit 1: \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\" heading ---\n$headingRange = $d.Content\n$headingFound = $headingRange.Find.Execute(\"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\")\nif (-not $headingFound) {\n    throw 'Could not find heading text \"NG\u01af\u1edcI \u0110\u1ea0I DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P\"'\n}\n\n$headingBody = '<w:p><w:pPr><w:spacing w:before=\"120\" w:after=\"120\"/><w:jc w:val=\"center\"/><w:rPr><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr><w:t>NG\u01af\u1edcI \u0110\u1ea0I</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr><w:t xml:space=\"preserve\"> DI\u1ec6N THEO PH\u00c1P LU\u1eacT C\u1ee6A DOANH NGHI\u1ec6P</w:t></w:r></w:p>'\n$headingRange.InsertXML($pkgOpen + $headingBody + $pkgClose)\n\n# --- Edit 2: \"{dissolution_approve_representative | upper}\" placeholder ---\n$placeholderRange = $d.Content\n$placeholderFound = $placeholderRange.Find.Execute(\"{dissolution_approve_representative | upper}\")\nif (-not $placeholderFound) {\n    throw 'Could not find placeholder text \"{dissolution_approve_representative | upper}\"'\n}\n\n$placeholderBody = '<w:p><w:pPr><w:spacing w:before=\"120\" w:after=\"120\"/><w:jc w:val=\"center\"/><w:rPr><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr><w:t xml:space=\"preserve\">dissolution_approve_org_person </w:t></w:r><w:r><w:rPr><w:sz w:val=\"26\"/><w:szCs w:val=\"26\"/><w:lang w:val=\"pt-BR\"/></w:rPr><w:t>| upper}</w:t></w:r></w:p>'\n$placeholderRange.InsertXML($pkgOpen + $placeholderBody + $pkgClose)\n"}
